$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.764.22'
$ws.Range("E2").Value = '  -0.79%  '
$ws.Range("D3").Value = '2.671.48'
$ws.Range("E3").Value = '  -0.47%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.87'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.43%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  +4.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.129'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.38%  '
$ws.Range("E10").Value = '  -1.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.84'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.71%  '
$ws.Range("E12").Value = '  -0.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.58'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000199'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.78%  '
$ws.Range("D15").Value = '3.151.65'
$ws.Range("E15").Value = '  -0.61%  '
$ws.Range("D16").Value = '65.585.76'
$ws.Range("E16").Value = '  -0.77%  '
$ws.Range("D17").Value = '2.682.67'
$ws.Range("E17").Value = '  -0.85%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.83'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.57%  '
$ws.Range("E19").Value = '  -1.68%  '
$ws.Range("E20").Value = '  +2.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '352.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.94%  '
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.90'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000110'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.84'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.87%  '
$ws.Range("E26").Value = '  -3.42%  '
$ws.Range("E27").Value = '  -3.20%  '
$ws.Range("E28").Value = '  -4.46%  '
$ws.Range("E29").Value = '  -1.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.16'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.68%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '530.11'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.50%  '
$ws.Range("E33").Value = '  -2.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.56'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.69%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.52'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.14%  '
$ws.Range("E36").Value = '  -2.25%  '
$ws.Range("E37").Value = '  -1.83%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '159.59'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.08%  '
$ws.Range("E39").Value = '  -0.09%  '
$ws.Range("E40").Value = '  -4.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '42.75'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '165.44'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.17%  '
$ws.Range("E44").Value = '  -3.25%  '
$ws.Range("E45").Value = '  -1.59%  '
$ws.Range("E46").Value = '  -1.08%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.17'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.27%  '
$ws.Range("E48").Value = '  -1.98%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.644'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.101'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.29'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.28%  '
